# Update latest output (run 32)
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Schedule": refresh the pump schedule (3 rows instead of 4)
# ---------------------------------------------------------------------------
$sched = $wb.Worksheets.Item("Schedule")

# Row 2 gets new schedule-block values
$sched.Range("A2").Value = 46038.16666666666
$sched.Range("B2").Value = 46038.66666666666
$sched.Range("C2").Value = 12
$sched.Range("D2").Value = 45.36
$sched.Range("E2").Value = 1126.019388
$sched.Range("F2").Value = 24.82406058201058

# Row 3 absorbs the former row-4 time block with recomputed cost figures
$sched.Range("A3").Value = 46038.83333333334
$sched.Range("B3").Value = 46039
$sched.Range("C3").Value = 4
$sched.Range("D3").Value = 15.12
$sched.Range("E3").Value = 337.0658655
$sched.Range("F3").Value = 22.29271597222222

# The old row 4 no longer exists -> drop it (also shrinks the dimension to A1:F3)
$sched.Rows(4).Delete()

# ---------------------------------------------------------------------------
# Sheet "Detailed": refresh prices / pump status for the latest optimisation run
# ---------------------------------------------------------------------------
$det = $wb.Worksheets.Item("Detailed")

# Pump status flips OFF for the early-morning block that used to be ON
$det.Range("E3").Value = "OFF"
$det.Range("E4").Value = "OFF"
$det.Range("E5").Value = "OFF"
$det.Range("E6").Value = "OFF"
$det.Range("E7").Value = "OFF"
$det.Range("E8").Value = "OFF"
$det.Range("E9").Value = "OFF"

# ... and flips ON for the mid-morning block that used to be OFF
$det.Range("E11").Value = "ON"
$det.Range("E12").Value = "ON"

$det.Range("B13").Value = 85.95
$det.Range("E13").Value = "ON"

$det.Range("B14").Value = 79.95
$det.Range("E14").Value = "ON"

$det.Range("C15").Value = "historical"
$det.Range("E15").Value = "ON"

$det.Range("C16").Value = "historical"
$det.Range("E16").Value = "ON"

$det.Range("B17").Value = 50.38708
$det.Range("E17").Value = "ON"

$det.Range("B18").Value = 56.98

$det.Range("B19").Value = 55.25227

$det.Range("B21").Value = 36.05916

$det.Range("B24").Value = 36.06028

$det.Range("B32").Value = 35.85034

$det.Range("B33").Value = 24.46863

$det.Range("B34").Value = 19.24233

$det.Range("B35").Value = 10.31216

$det.Range("B36").Value = -1.5001

$det.Range("B38").Value = -3.17664

$det.Range("B39").Value = -2.7582

$det.Range("B40").Value = 0.0113

$det.Range("B43").Value = 29.85322

$det.Range("B44").Value = 0.84406

$det.Range("B45").Value = 57.04922

$det.Range("B47").Value = 57.04922

$det.Range("B48").Value = 57.03042
